# New test results for sha2-384 and sha2-512:
# - Row 4 (SHA_512) gets refreshed benchmark numbers.
# - Row 5's label is changed from SHA_256 to the new SHA_384 entry, with its
#   own fresh benchmark numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: SHA_512 (values only update, label stays the same) ---
$ws.Range("C4").Value = 14
$ws.Range("D4").Value = 2005
$ws.Range("E4").Value = 1809
$ws.Range("F4").Value = 1566

# --- Row 5: relabel to SHA_384 and refresh all benchmark values ---
$ws.Range("A5").Value = "SHA_384"
$ws.Range("C5").Value = 14
$ws.Range("D5").Value = 1758
$ws.Range("E5").Value = 1524
$ws.Range("F5").Value = 1139
$ws.Range("H5").Value = 1127
$ws.Range("I5").Value = 945
$ws.Range("J5").Value = 997

# --- Selection / view tidy-up to match the reopened workbook state ---
$ws.Range("A7").Select()
